# Add a new "10-nov" data column (CJ) to the sheet, right after the
# existing "9-nov" column (CI), with one new value per family row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: copy the text-style formatting used by the other date headers.
$ws.Cells.Item(1, 88).Value = "10-nov"
$ws.Cells.Item(1, 88).NumberFormat = $ws.Cells.Item(1, 87).NumberFormat

# New data values for each family (rows 2..11) in the new column.
$newValues = @(8, 12, 7, 11, 7, 6, 15, 15, 8, 0)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 88)
    $cell.Value = $newValues[$i]
    # Match the centered-integer formatting used by the rest of the data cells.
    $cell.HorizontalAlignment = $ws.Cells.Item($row, 87).HorizontalAlignment
    $cell.NumberFormat = $ws.Cells.Item($row, 87).NumberFormat
}

# Reflect where editing left off, as in the original edit.
$ws.Range("CJ10").Select()
